$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.738.67'
$ws.Range('E2').Value = '  +0.66%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.965.54'
$ws.Range('E3').Value = '  +1.11%  '

# Row 4
$ws.Range('E4').Value = '  -0.12%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.98'
$ws.Range('E5').Value = '  +0.55%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.620'
$ws.Range('E6').Value = '  +0.83%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.00'
$ws.Range('E7').Value = '  +1.90%  '

# Row 8
$ws.Range('E8').Value = '  -0.03%  '

# Row 9
$ws.Range('E9').Value = '  +1.97%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0816'
$ws.Range('E10').Value = '  -2.77%  '

# Row 11
$ws.Range('E11').Value = '  -0.33%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.33'
$ws.Range('E12').Value = '  +3.05%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.256.81'
$ws.Range('E13').Value = '  +1.11%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.828'
$ws.Range('E14').Value = '  +0.42%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '13.75'
$ws.Range('E15').Value = '  +1.09%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.27'
$ws.Range('E16').Value = '  +0.27%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.961.59'
$ws.Range('E17').Value = '  +1.48%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.669.31'
$ws.Range('E18').Value = '  +0.61%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.82'
$ws.Range('E19').Value = '  +0.18%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0861'
$ws.Range('E20').Value = '  -1.05%  '

# Row 21
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '229.17'
$ws.Range('E21').Value = '  -0.33%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.08'
$ws.Range('E22').Value = '  +1.40%  '

# Row 23
$ws.Range('E23').Value = '  -0.08%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.44'
$ws.Range('E24').Value = '  -0.33%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.37'
$ws.Range('E25').Value = '  +2.96%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.32'
$ws.Range('E26').Value = '  +0.42%  '

# Row 27
$ws.Range('E27').Value = '  +14.18%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.52'
$ws.Range('E28').Value = '  -1.07%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.39'
$ws.Range('E29').Value = '  -0.20%  '

# Row 30
$ws.Range('E30').Value = '  +1.29%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.13'
$ws.Range('E31').Value = '  -1.70%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.71'
$ws.Range('E32').Value = '  +0.76%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0619'
$ws.Range('E33').Value = '  -1.69%  '

# Row 34
$ws.Range('E34').Value = '  +0.47%  '

# Row 35
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.15%  '

# Row 36
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.27'
$ws.Range('E36').Value = '  +5.57%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.09'
$ws.Range('E37').Value = '  -2.81%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.39'
$ws.Range('E38').Value = '  +11.06%  '

# Row 39
$ws.Range('E39').Value = '  -0.14%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.100'
$ws.Range('E40').Value = '  +3.10%  '

# Row 41
$ws.Range('E41').Value = '  -2.26%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0212'
$ws.Range('E42').Value = '  +1.53%  '

# Row 43
$ws.Range('E43').Value = '  -0.81%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.11'
$ws.Range('E44').Value = '  +0.28%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.368.15'
$ws.Range('E45').Value = '  +0.92%  '

# Row 46
$ws.Range('E46').Value = '  +0.76%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '87.89'
$ws.Range('E47').Value = '  +0.17%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.15'
$ws.Range('E48').Value = '  -0.16%  '

# Row 49
$ws.Range('E49').Value = '  +0.67%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.148.00'
$ws.Range('E50').Value = '  +1.23%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.75'
$ws.Range('E51').Value = '  -3.28%  '
